$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1196
$ws1.Range("F5").Value = 70
$ws1.Range("F8").Value = 1935
$ws1.Range("F10").Value = 558
$ws1.Range("F12").Value = 87
$ws1.Range("F14").Value = 747
$ws1.Range("F15").Value = 532
$ws1.Range("F16").Value = 931
$ws1.Range("F17").Value = 86259
$ws1.Range("F18").Value = 3
$ws1.Range("F21").Value = 39463
$ws1.Range("F22").Value = 637
$ws1.Range("F23").Value = 51
$ws1.Range("F27").Value = 1104
$ws1.Range("F28").Value = 41
$ws1.Range("F29").Value = 364
$ws1.Range("F31").Value = 776
$ws1.Range("F32").Value = 72
$ws1.Range("F34").Value = 1312
$ws1.Range("F35").Value = 5617
$ws1.Range("F37").Value = 511
$ws1.Range("F43").Value = 552

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F36").Value = 62
$ws2.Range("F44").Value = 364
$ws2.Range("F47").Value = 6

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 612
$ws3.Range("F7").Value = 268
$ws3.Range("F8").Value = 132

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 612
$ws4.Range("F6").Value = 1196
$ws4.Range("F9").Value = 70
$ws4.Range("F11").Value = 268
$ws4.Range("F12").Value = 132
$ws4.Range("F13").Value = 132
$ws4.Range("F14").Value = 1935
$ws4.Range("F16").Value = 558
$ws4.Range("F18").Value = 87
$ws4.Range("F20").Value = 747
$ws4.Range("F23").Value = 533
$ws4.Range("F24").Value = 931
$ws4.Range("F26").Value = 86262
$ws4.Range("F27").Value = 3
$ws4.Range("F29").Value = 39463
$ws4.Range("F30").Value = 637
$ws4.Range("F31").Value = 51
$ws4.Range("F33").Value = 1104
$ws4.Range("F34").Value = 41
$ws4.Range("F37").Value = 364
$ws4.Range("F38").Value = 776
$ws4.Range("F39").Value = 72
$ws4.Range("F40").Value = 1312
$ws4.Range("F41").Value = 5617
$ws4.Range("F44").Value = 511
$ws4.Range("F49").Value = 552
$ws4.Range("F53").Value = 364

$wb.Save()
